$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 17: "NEW_MESSAGE-1" -> "COMPOSEMESSAGE" (Runmode stays "N")
$ws.Cells.Item(17, 1).Value = "COMPOSEMESSAGE"

# Row 18: "COMPOSEMESSAGE" -> "SIGNIN" (Runmode stays "Y")
$ws.Cells.Item(18, 1).Value = "SIGNIN"
$ws.Cells.Item(18, 3).Value = "Y"

# Update the active selection to H13 (as recorded in the saved view state)
$null = $ws.Range("H13").Select()

# Force a full recalculation on load of the workbook
$wb.ForceFullCalculation = $true

# Restore/request the window size recorded in the saved view state
$win = $excel.ActiveWindow
if ($win) {
    $win.Height = 13240
}
